# Rename the 196-numbered Manufacturer/AuthorisedRep test users to 371, and
# flip the "ignore" flag for the AT / NU row groups (AT: no -> yes,
# NU: yes -> no), per the commit:
#   Created 4 new users: Manufacturer371_AT, Manufacturer371_NU,
#   AuthorisedRep371_AT, AuthorisedRep371_NU

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# AT rows (3-4: Manufacturer/AuthorisedRep usernames renumbered) + ignore flag -> yes
$ws.Range("C2").Value = "yes"

$ws.Range("A3").Value = "Manufacturer371_AT"
$ws.Range("C3").Value = "yes"

$ws.Range("A4").Value = "AuthorisedRep371_AT"
$ws.Range("C4").Value = "yes"

# NU rows (6-7: Manufacturer/AuthorisedRep usernames renumbered) + ignore flag -> no
$ws.Range("C5").Value = "no"

$ws.Range("A6").Value = "Manufacturer371_NU"
$ws.Range("C6").Value = "no"

$ws.Range("A7").Value = "AuthorisedRep371_NU"
$ws.Range("C7").Value = "no"

# Move the saved selection from B9 to A8
$ws.Range("A8").Select()
